$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2163
$ws1.Range("F5").Value = 11401
$ws1.Range("F7").Value = 319
$ws1.Range("F9").Value = 11349
$ws1.Range("F10").Value = 460
$ws1.Range("F12").Value = 71
$ws1.Range("F13").Value = 1744
$ws1.Range("F14").Value = 5660
$ws1.Range("F15").Value = 106

# Sheet "全部类型" (sheet4.xml) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2163
$ws4.Range("F7").Value = 11401
$ws4.Range("F9").Value = 319
$ws4.Range("F11").Value = 11349
$ws4.Range("F12").Value = 460
$ws4.Range("F14").Value = 71
$ws4.Range("F15").Value = 1744
$ws4.Range("F17").Value = 5660
$ws4.Range("F18").Value = 106
